{"js": "// Replace multiplication-problem text in each table cell according to the\n// mapping derived from the diff. Each \"from\" string is unique in the\n// document, so a simple search/replace per pair is safe and unambiguous.\nconst replacements = [\n  [\"82\u00d761=\", \"76\u00d781=\"],\n  [\"16\u00d791=\", \"98\u00d716=\"],\n  [\"49\u00d740=\", \"35\u00d743=\"],\n  [\"53\u00d739=\", \"95\u00d715=\"],\n  [\"62\u00d738=\", \"96\u00d739=\"],\n  [\"81\u00d719=\", \"54\u00d778=\"],\n  [\"49\u00d738=\", \"92\u00d769=\"],\n  [\"12\u00d745=\", \"36\u00d727=\"],\n  [\"81\u00d758=\", \"14\u00d724=\"],\n  [\"25\u00d779=\", \"23\u00d774=\"],\n  [\"59\u00d737=\", \"94\u00d722=\"],\n  [\"83\u00d788=\", \"98\u00d779=\"],\n  [\"30\u00d787=\", \"50\u00d776=\"],\n  [\"25\u00d735=\", \"92\u00d721=\"],\n  [\"92\u00d798=\", \"13\u00d793=\"],\n  [\"76\u00d722=\", \"50\u00d750=\"],\n  [\"78\u00d746=\", \"78\u00d747=\"],\n  [\"92\u00d744=\", \"81\u00d732=\"],\n  [\"85\u00d766=\", \"91\u00d774=\"],\n  [\"98\u00d777=\", \"82\u00d782=\"],\n  [\"20\u00d718=\", \"86\u00d768=\"],\n  [\"62\u00d771=\", \"59\u00d773=\"],\n  [\"40\u00d775=\", \"40\u00d728=\"],\n  [\"22\u00d793=\", \"96\u00d729=\"],\n  [\"69\u00d752=\", \"38\u00d777=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [from, to] of replacements) {\n  const results = body.search(from, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(to, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace multiplication-problem text in each table cell according to the\n# mapping derived from the diff. Each \"from\" string is unique in the\n# document, so Find/Replace (wdReplaceAll) per pair is safe and unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"82\u00d761=\", \"76\u00d781=\"),\n    @(\"16\u00d791=\", \"98\u00d716=\"),\n    @(\"49\u00d740=\", \"35\u00d743=\"),\n    @(\"53\u00d739=\", \"95\u00d715=\"),\n    @(\"62\u00d738=\", \"96\u00d739=\"),\n    @(\"81\u00d719=\", \"54\u00d778=\"),\n    @(\"49\u00d738=\", \"92\u00d769=\"),\n    @(\"12\u00d745=\", \"36\u00d727=\"),\n    @(\"81\u00d758=\", \"14\u00d724=\"),\n    @(\"25\u00d779=\", \"23\u00d774=\"),\n    @(\"59\u00d737=\", \"94\u00d722=\"),\n    @(\"83\u00d788=\", \"98\u00d779=\"),\n    @(\"30\u00d787=\", \"50\u00d776=\"),\n    @(\"25\u00d735=\", \"92\u00d721=\"),\n    @(\"92\u00d798=\", \"13\u00d793=\"),\n    @(\"76\u00d722=\", \"50\u00d750=\"),\n    @(\"78\u00d746=\", \"78\u00d747=\"),\n    @(\"92\u00d744=\", \"81\u00d732=\"),\n    @(\"85\u00d766=\", \"91\u00d774=\"),\n    @(\"98\u00d777=\", \"82\u00d782=\"),\n    @(\"20\u00d718=\", \"86\u00d768=\"),\n    @(\"62\u00d771=\", \"59\u00d773=\"),\n    @(\"40\u00d775=\", \"40\u00d728=\"),\n    @(\"22\u00d793=\", \"96\u00d729=\"),\n    @(\"69\u00d752=\", \"38\u00d777=\")\n)\n\nforeach ($pair in $replacements) {\n    $from = $pair[0]\n    $to = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $from\n    $find.Replacement.Text = $to\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute(\n        $find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards,\n        $null, $null, $find.Forward, $find.Wrap, $null,\n        $find.Replacement.Text, 2  # wdReplaceAll\n    )\n}\n"}
